$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contact person changed
$ws.Range("B7").Value = "Mambetaliev T.A."

# 2. Data reporter -> Organization renamed
$ws.Range("B6").Value = "The National Statistical Committee of the Kyrgyz Republic (Department of Digital Development and Sustainable Development Statistics)"

# Contact person's phone changed
$ws.Range("B9").Value = "(0312) 62 56 07"

# Organization website changed
$ws.Range("B10").Value = "www.stat.gov.kg"

# 1. Indicator information -> Indicator text update (12.4.2 wording revision)
$ws.Range("B4").Value = "12.4.2 (a) Hazardous waste generated per capita; and (b) proportion of hazardous waste treated, by type of treatment"

# Update the active selection to match the new focus cell
$ws.Range("B4").Select()
